$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 rows before the (currently) second summary table (old row 13),
# pushing it down to rows 17-18, and making room for the new commentary
# lines at rows 13-16.
$ws.Rows.Item(13).Resize(4).Insert()

$ws.Range("B13").Value = "|"
$ws.Range("B14").Value = "Remove one outlier in HDL"
$ws.Range("B15").Value = "|"
$ws.Range("B16").Value = "1776 rows remained."

# New pairwise-test / residual-check output block, rows 20-24, monospaced.
$ws.Range("B20").Value = "          F   M"
$ws.Range("B21").Value = "  DivNA 467 419"
$ws.Range("B22").Value = "  Div0  311 252"
$ws.Range("B23").Value = "  Div1   88  81"
$ws.Range("B24").Value = "  Div2   82  76"

# Apply the monospaced font via a transient named style so the engine
# resolves straight to the final font (avoids leaving an intermediate,
# unused "de-schemed default font" entry in the style table that a
# direct `.Font.Name =` assignment on the range would otherwise leave
# behind).
$consolasStyle = $wb.Styles.Add("ConsolasTemp")
$consolasStyle.Font.Name = "Consolas"
$ws.Range("B20:B24").Style = "ConsolasTemp"
$wb.Styles.Item("ConsolasTemp").Delete()

$ws.Range("D26").Select()
